$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column (C) for rows 2-8
# from serial 45170 (2023-09-01) to serial 45174 (2023-09-05)
foreach ($row in 2..8) {
    $ws.Cells.Item($row, 3).Value = 45174
}
